$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Portfolio")
$ws.Range("A12").Value = "Test"
